$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- old Row 3 values
$ws.Range("D2").Value = 44568
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 16000
$ws.Range("M2").Value = 15500
$ws.Range("P2").Value = 861

# Row 3 <- old Row 4 values
$ws.Range("D3").Value = 44547
$ws.Range("J3").Value = 200
$ws.Range("K3").Value = 13000
$ws.Range("L3").Value = 14000
$ws.Range("M3").Value = 13500
$ws.Range("P3").Value = 750

# Row 4 <- old Row 2 values
$ws.Range("D4").Value = 44557
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 13000
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = 13500
$ws.Range("P4").Value = 750
